# =====================================================================
# Create Employee Screen Integration Basic
# =====================================================================
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the two existing sheets.
# ---------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item(1)
$wsTasks.Name = "Tasks"

$wsBugs = $wb.Worksheets.Item(2)
$wsBugs.Name = "Bugs"

# ---------------------------------------------------------------------
# 2. Add the two new sheets so the final tab order is:
#    Tasks, Bugs, Test data backup, create employee screen
#    (create employee screen is added first so it gets sheetId 3,
#    Test data backup is added second so it gets sheetId 4, then it is
#    moved in front of "create employee screen").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCES = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCES.Name = "create employee screen"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTDB = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsTDB.Name = "Test data backup"

# Move "Test data backup" so it sits before "create employee screen".
$wsTDB.Move($wsCES)

# Re-fetch sheet references by name now that the tab order changed.
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsBugs  = $wb.Worksheets.Item("Bugs")
$wsTDB   = $wb.Worksheets.Item("Test data backup")
$wsCES   = $wb.Worksheets.Item("create employee screen")

# ---------------------------------------------------------------------
# 3. Tasks sheet: two new rows of bug-fix follow-up work.
# ---------------------------------------------------------------------
$wsTasks.Range("B22").Value = 45276
$wsTasks.Range("C22").Value = 45276

$wsTasks.Range("B23").Value = 45276
$wsTasks.Range("C23").Value = 45277
$wsTasks.Range("D23").Value = "open bugs = fixed"

$wsTasks.Range("B22:C23").NumberFormat = "mm-dd-yy"

$wsTasks.Columns("B:C").ColumnWidth = 9.6

# ---------------------------------------------------------------------
# 4. Bugs sheet: one new bug row.
# ---------------------------------------------------------------------
$wsBugs.Range("A3").Value = "in save api the req.body is printing an empty obj"
$wsBugs.Range("B3").Value = "fixed"

# ---------------------------------------------------------------------
# 5. Test data backup sheet: sample employee used while testing the
#    create-employee screen.
# ---------------------------------------------------------------------
$wsTDB.Range("A1").Value = "A-300"
$wsTDB.Range("A2").Value = "Anish"
$wsTDB.Range("A3").Value = 9447757072
$wsTDB.Range("A4").Value = "Anish Bhavan"

$wsTDB.Columns("A").ColumnWidth = 39.5

# ---------------------------------------------------------------------
# 6. create employee screen sheet: the actual employee-creation form
#    layout (section headers + field labels).
# ---------------------------------------------------------------------

# -- Row 1: top title ---------------------------------------------------
$wsCES.Range("A1").Value = "Sl.No"

# -- Row 3: table header --------------------------------------------
$wsCES.Range("B3").Value = "Employee code"
$wsCES.Range("C3").Value = "Employee Name"

# -- Row 4: section headers ------------------------------------------
$wsCES.Range("A4").Value = "Contact Information"
$wsCES.Range("K4").Value = "Employment Information"
$wsCES.Range("P4").Value = "Personal Information"

# -- Row 5: field labels ----------------------------------------------
$wsCES.Range("B5").Value = "Address"
$wsCES.Range("C5").Value = "Mobile Number"
$wsCES.Range("K5").Value = "Designation"
$wsCES.Range("L5").Value = "Date of joining"
$wsCES.Range("M5").Value = "Status"
$wsCES.Range("P5").Value = "Name of father / husband"
$wsCES.Range("Q5").Value = "Gender"
$wsCES.Range("R5").Value = "Date of Birth"
$wsCES.Range("S5").Value = "Marital Status"

# -- Row 7: more employment fields -------------------------------------
$wsCES.Range("K7").Value = "Last Working Date"
$wsCES.Range("L7").Value = "Date of Releaving"
$wsCES.Range("M7").Value = "Remarks"

# -- Row 10/11: bank information ---------------------------------------
$wsCES.Range("A10").Value = "Bank Information"
$wsCES.Range("B11").Value = "Bank Name"
$wsCES.Range("C11").Value = "IFSC Code"
$wsCES.Range("D11").Value = "Bank Account Number"

# -- Row 12/13: identity information ------------------------------------
$wsCES.Range("A12").Value = "Identity Information"
$wsCES.Range("B13").Value = "Pan Number"
$wsCES.Range("C13").Value = "Aadhar Number"
$wsCES.Range("D13").Value = "UAN EPF"
$wsCES.Range("E13").Value = "EPF No"
$wsCES.Range("F13").Value = "ESI Number"
$wsCES.Range("G13").Value = "Welfare Fund Number"

# -- Row 15+: nominees ---------------------------------------------------
$wsCES.Range("A15").Value = "Nominees "
$wsCES.Range("B16").Value = "EPF Nominee"
$wsCES.Range("C16").Value = "EPF Nominee Relation"
$wsCES.Range("B18").Value = "ESI Nominee"
$wsCES.Range("C18").Value = "ESI Nominee Relation"
$wsCES.Range("B20").Value = "GPAIP Nominee"
$wsCES.Range("C20").Value = "GPAIP Nominee Relation"
$wsCES.Range("B22").Value = "Gratuity Nominee"
$wsCES.Range("C22").Value = "Gratuity Nominee Relation"

# Row heights for the wrapped header rows.
$wsCES.Rows("3").RowHeight = 24
$wsCES.Rows("5").RowHeight = 36
$wsCES.Rows("7").RowHeight = 36
$wsCES.Rows("11").RowHeight = 36
$wsCES.Rows("13").RowHeight = 36
$wsCES.Rows("16").RowHeight = 36
$wsCES.Rows("18").RowHeight = 36
$wsCES.Rows("20").RowHeight = 36
$wsCES.Rows("22").RowHeight = 36

# -----------------------------------------------------------------
# Formatting: the "Ahalf" (9pt) styled header/label cells with the
# themed fill + thin border used throughout the form.
# -----------------------------------------------------------------
$ahamCenterWrap = $excel.Union($wsCES.Range("A1"), $wsCES.Range("C5"))
$ahamCenterWrap = $excel.Union($ahamCenterWrap, $wsCES.Range("K5"))
$ahamCenterWrap = $excel.Union($ahamCenterWrap, $wsCES.Range("C13"))
$ahamCenterWrap = $excel.Union($ahamCenterWrap, $wsCES.Range("D13"))

$ahamLeftWrap = $excel.Union($wsCES.Range("B3"), $wsCES.Range("C3"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("M5"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("P5"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("S5"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("B5"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("B11"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("C11"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("D11"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("B13"))
$ahamLeftWrap = $excel.Union($ahamLeftWrap, $wsCES.Range("E13"))

$ahamVertOnly = $wsCES.Range("Q5")

$ahamDateLeft = $wsCES.Range("R5")

$ahamDateCenterWrap = $wsCES.Range("L5")

$arialCenterWrap = $wsCES.Range("F13")

$arialLeftWrap = $excel.Union($wsCES.Range("K7"), $wsCES.Range("L7"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("M7"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("G13"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("B16"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("C16"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("B18"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("C18"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("B20"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("C20"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("B22"))
$arialLeftWrap = $excel.Union($arialLeftWrap, $wsCES.Range("C22"))

$allAham = $excel.Union($ahamCenterWrap, $ahamLeftWrap)
$allAham = $excel.Union($allAham, $ahamVertOnly)
$allAham = $excel.Union($allAham, $ahamDateLeft)
$allAham = $excel.Union($allAham, $ahamDateCenterWrap)

$allArial = $excel.Union($arialCenterWrap, $arialLeftWrap)

$allLabels = $excel.Union($allAham, $allArial)

# Shared look: themed fill + thin border for every label cell.
$allLabels.Interior.Pattern = -4124
$allLabels.Interior.ThemeColor = 6
$allLabels.Interior.TintAndShade = 0.59999389629810485
$allLabels.Borders.LineStyle = 1
$allLabels.Borders.Weight = 2

$allAham.Font.Name = "A\ham"
$allAham.Font.Size = 9

$allArial.Font.Name = "Arial"
$allArial.Font.Size = 9

# Alignment per group.
$ahamCenterWrap.HorizontalAlignment = -4108
$ahamCenterWrap.VerticalAlignment = -4108
$ahamCenterWrap.WrapText = $true

$ahamLeftWrap.HorizontalAlignment = -4131
$ahamLeftWrap.VerticalAlignment = -4108
$ahamLeftWrap.WrapText = $true

$ahamVertOnly.VerticalAlignment = -4108

$ahamDateLeft.HorizontalAlignment = -4131
$ahamDateLeft.VerticalAlignment = -4108
$ahamDateLeft.NumberFormat = "mm-dd-yy"

$ahamDateCenterWrap.HorizontalAlignment = -4108
$ahamDateCenterWrap.VerticalAlignment = -4108
$ahamDateCenterWrap.WrapText = $true
$ahamDateCenterWrap.NumberFormat = "mm-dd-yy"

$arialCenterWrap.HorizontalAlignment = -4108
$arialCenterWrap.VerticalAlignment = -4108
$arialCenterWrap.WrapText = $true

$arialLeftWrap.HorizontalAlignment = -4131
$arialLeftWrap.VerticalAlignment = -4108
$arialLeftWrap.WrapText = $true

$wsCES.Columns("G").ColumnWidth = 10.6

$wsCES.Range("A12").Select()
$wsCES.Application.ActiveWindow.ScrollRow = 10

# ---------------------------------------------------------------------
# 7. Final selections per-sheet, and re-activate "Tasks" so it is the
#    tab shown when the workbook is opened (matches the original file).
# ---------------------------------------------------------------------
$wsBugs.Range("B4").Select()
$wsTDB.Range("A6").Select()
$wsCES.Range("A12").Select()

$wsTasks.Range("D23").Select()
$wsTasks.Activate()

Write-Output "done"
